$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Add back fastq rows that were found in LTS but had no metadata files.
# Columns: A harvestDate, B harvester, C bioSampleNumber, D rnaDate,
#          E rnaPreparer, F rnaSampleNumber, G rnaPrepMethod, H roboticRNAPrep
# ---------------------------------------------------------------------------

# Row 35 mirrors the pattern of row 2 (harvestDate/rnaDate = 08.09.18),
# except its rnaDate becomes the new "08.13.18" value.
$ws.Range("A2:H2").Copy($ws.Range("A35:H35"))
$ws.Range("C35").Value = 35
$ws.Range("F35").Value = 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Characters(4, 2).Text = "13"
$ws.Range("D35").Style = $ws.Range("B2").Style

# Rows 36-37 mirror the pattern of row 3 (harvestDate/rnaDate = 10.15.18),
# except their rnaDate becomes the new "10.16.18" value.
$ws.Range("A3:H3").Copy($ws.Range("A36:H36"))
$ws.Range("C36").Value = 36
$ws.Range("F36").Value = 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Characters(5, 1).Text = "6"
$ws.Range("D36").Style = $ws.Range("B2").Style

$ws.Range("A3:H3").Copy($ws.Range("A37:H37"))
$ws.Range("C37").Value = 37
$ws.Range("F37").Value = 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Characters(5, 1).Text = "6"
$ws.Range("D37").Style = $ws.Range("B2").Style

# (G35:G37 / H35:H37 already inherited the correct "TRIzol"/"False" styling
# -- style indices 2/3 -- from the row 2/row 3 copy above; re-touching
# .Style here would only risk Excel re-resolving it away.)

# Row 38: a trailing, mostly-empty row with only H38 carrying the text
# number format (no value).
$ws.Range("H38").NumberFormat = "@"

# Move the active selection to match where the author ended up (E38).
$ws.Range("E38").Select()
